# Swap the "step 2" content between TC2 and TC3 test cases.
# TC2's second step currently describes the "cancel a diária" action; TC3's
# second step currently describes the "filter/search" action. This change
# exchanges those two step descriptions so that TC2 now contains the
# filter/search step and TC3 contains the cancel step (labels TC2/TC3 stay
# where they are).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cancelSteps   = $ws.Range("B20").Value2
$cancelResults = $ws.Range("D20").Value2

$searchSteps   = $ws.Range("B28").Value2
$searchResults = $ws.Range("D28").Value2

$ws.Range("B20").Value2 = $searchSteps
$ws.Range("D20").Value2 = $searchResults

$ws.Range("B28").Value2 = $cancelSteps
$ws.Range("D28").Value2 = $cancelResults
